$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1904912.1
$ws.Range("I19").Value = 2666776.5
$ws.Range("J19").Value = 251
$ws.Range("K19").Value = 2666776.5
$ws.Range("L19").Value = 251
$ws.Range("M19").Value = -2666601.5
$ws.Range("N19").Value = -601
$ws.Range("H62").Value = 2942.625
$ws.Range("I62").Value = 2257
$ws.Range("K62").Value = 2257
$ws.Range("M62").Value = -1633
$ws.Range("H64").Value = 3100
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("H65").Value = 2942.625
$ws.Range("I65").Value = 2257
$ws.Range("K65").Value = 11285
$ws.Range("M65").Value = -8165
$ws.Range("H67").Value = 3100
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("H86").Value = 1610.7059
$ws.Range("I86").Value = 391
$ws.Range("J86").Value = 5574.75
$ws.Range("K86").Value = 391
$ws.Range("L86").Value = 5574.75
$ws.Range("M86").Value = 732
$ws.Range("N86").Value = -7820.75
$ws.Range("H89").Value = 1610.7059
$ws.Range("I89").Value = 391
$ws.Range("J89").Value = 5574.75
$ws.Range("K89").Value = 1955
$ws.Range("L89").Value = 27873.75
$ws.Range("M89").Value = 3661
$ws.Range("N89").Value = -39105.75
$ws.Range("H111").Value = 1500
$ws.Range("I111").Value = 1500
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 4500
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -1433
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 483348.34
$ws.Range("J116").Value = 11175.417
$ws.Range("L116").Value = 11175.417
$ws.Range("N116").Value = -18059.417
$ws.Range("H118").Value = 1945.3846
$ws.Range("I118").Value = 1830
$ws.Range("J118").Value = 1980
$ws.Range("K118").Value = 5490
$ws.Range("L118").Value = 5940
$ws.Range("M118").Value = -3833
$ws.Range("N118").Value = -9254
$ws.Range("H129").Value = 942.6512
$ws.Range("J129").Value = 969.12195
$ws.Range("L129").Value = 2907.36585
$ws.Range("N129").Value = -12907.36585
$ws.Range("H138").Value = 2172.36
$ws.Range("I138").Value = 989.93616
$ws.Range("J138").Value = 3220.9246
$ws.Range("K138").Value = 2969.80848
$ws.Range("L138").Value = 9662.773799999999
$ws.Range("M138").Value = 2170.19152
$ws.Range("N138").Value = -19942.7738
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 513.7292
$ws.Range("I2").Value = 487.64102
$ws.Range("K2").Value = 487.64102
$ws.Range("M2").Value = -374.64102
$ws.Range("H32").Value = 8441.111999999999
$ws.Range("I32").Value = 5410.379
$ws.Range("J32").Value = 14692
$ws.Range("K32").Value = 5410.379
$ws.Range("L32").Value = 14692
$ws.Range("M32").Value = -5123.379
$ws.Range("N32").Value = -15266
$ws.Range("H45").Value = 2327.2
$ws.Range("I45").Value = 1156
$ws.Range("J45").Value = 2620
$ws.Range("K45").Value = 1156
$ws.Range("L45").Value = 2620
$ws.Range("M45").Value = -779
$ws.Range("N45").Value = -3374
$ws.Range("H63").Value = 7292367
$ws.Range("I63").Value = 8659085
$ws.Range("J63").Value = 3202
$ws.Range("K63").Value = 8659085
$ws.Range("L63").Value = 3202
$ws.Range("M63").Value = -8658399
$ws.Range("N63").Value = -4574
$ws.Range("H66").Value = 7292367
$ws.Range("I66").Value = 8659085
$ws.Range("J66").Value = 3202
$ws.Range("K66").Value = 43295425
$ws.Range("L66").Value = 16010
$ws.Range("M66").Value = -43291993
$ws.Range("N66").Value = -22874
$ws.Range("H116").Value = 513.7292
$ws.Range("I116").Value = 487.64102
$ws.Range("K116").Value = 487.64102
$ws.Range("M116").Value = 1806.35898
$ws.Range("H122").Value = 2380.6428
$ws.Range("I122").Value = 1409.9231
$ws.Range("K122").Value = 4229.7693
$ws.Range("M122").Value = -1779.7693
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 513.7292
$ws.Range("I3").Value = 487.64102
$ws.Range("K3").Value = 487.64102
$ws.Range("M3").Value = -373.64102
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3547.25
$ws.Range("I105").Value = 10000
$ws.Range("K105").Value = 10000
$ws.Range("M105").Value = -8253
$ws.Range("H134").Value = 4585.8857
$ws.Range("I134").Value = 4890.88
$ws.Range("K134").Value = 14672.64
$ws.Range("M134").Value = -12137.64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 12130.889
$ws.Range("I68").Value = 956
$ws.Range("J68").Value = 26099.5
$ws.Range("K68").Value = 2868
$ws.Range("L68").Value = 78298.5
$ws.Range("M68").Value = -2057
$ws.Range("N68").Value = -79920.5
$ws.Range("H71").Value = 12130.889
$ws.Range("I71").Value = 956
$ws.Range("J71").Value = 26099.5
$ws.Range("K71").Value = 8604
$ws.Range("L71").Value = 234895.5
$ws.Range("M71").Value = -4548
$ws.Range("N71").Value = -243007.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 35716644
$ws.Range("I80").Value = 62501876
$ws.Range("K80").Value = 62501876
$ws.Range("M80").Value = -62500878
$ws.Range("H83").Value = 35716644
$ws.Range("I83").Value = 62501876
$ws.Range("K83").Value = 312509380
$ws.Range("M83").Value = -312504388
$ws.Range("H122").Value = 3174.8333
$ws.Range("I122").Value = 2773.125
$ws.Range("K122").Value = 8319.375
$ws.Range("M122").Value = -5869.375
$ws.Range("H126").Value = 3660.33
$ws.Range("I126").Value = 2719.8872
$ws.Range("J126").Value = 5194.737
$ws.Range("K126").Value = 8159.6616
$ws.Range("L126").Value = 15584.211
$ws.Range("M126").Value = -5689.6616
$ws.Range("N126").Value = -20524.211
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 925
$ws.Range("I68").Value = 840.1316
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 840.1316
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -91.13160000000005
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 925
$ws.Range("I71").Value = 840.1316
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 4200.658
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -456.6580000000004
$ws.Range("N71").Value = -17488
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 2790.946
$ws.Range("I136").Value = 1394.375
$ws.Range("J136").Value = 5369.231
$ws.Range("K136").Value = 4183.125
$ws.Range("L136").Value = 16107.693
$ws.Range("M136").Value = -1633.125
$ws.Range("N136").Value = -21207.693
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3506.0588
$ws.Range("I136").Value = 1158.0625
$ws.Range("J136").Value = 5593.1665
$ws.Range("K136").Value = 3474.1875
$ws.Range("L136").Value = 16779.4995
$ws.Range("M136").Value = -924.1875
$ws.Range("N136").Value = -21879.4995
